# Got a delay timer working to give items time to drop into tray before
# moving it.
#
# The "EA"/"EB" rows (18 & 19) in the wiring table used to document the
# two extra breadboard/Atmega pins ("29/30", "50/51") and note
# "PORTB:2 or 3" for the delay-timer wiring. Now that the delay timer is
# wired up and working, those rows just need a simple "5V" reference in
# the Breadboard Pin column (D), matching the other power rows, with the
# Atmega Register (E) and Notes (F) columns cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "5V"
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""

$ws.Range("D19").Value = "5V"
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""

# Leave the selection on the delay-timer row group (A17:A21) where the
# edits were made, with A17 as the active cell.
$ws.Range("A17:A21").Select()
